$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.820.33'
$ws.Range('E2').Value = '  +2.06%  '
$ws.Range('D3').Value = '3.089.00'
$ws.Range('E3').Value = '  +5.13%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.71'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '167.87'
$ws.Range('E6').Value = '  +5.75%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.086.02'
$ws.Range('E8').Value = '  +5.18%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.155'
$ws.Range('E11').Value = '  +3.29%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.483'
$ws.Range('E12').Value = '  +5.37%  '
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.39'
$ws.Range('E14').Value = '  +6.09%  '
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '3.602.45'
$ws.Range('E16').Value = '  +5.13%  '
$ws.Range('D17').Value = '66.822.51'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.20'
$ws.Range('E18').Value = '  +3.75%  '
$ws.Range('D19').Value = '3.090.19'
$ws.Range('E19').Value = '  +5.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.09'
$ws.Range('E20').Value = '  +3.14%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '467.00'
$ws.Range('E21').Value = '  +4.94%  '
$ws.Range('E22').Value = '  +3.71%  '
$ws.Range('E23').Value = '  +3.42%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '83.84'
$ws.Range('E24').Value = '  +1.99%  '
$ws.Range('E25').Value = '  +6.29%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '13.06'
$ws.Range('E26').Value = '  +7.95%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.12'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('E31').Value = '  +3.94%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0000102'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '28.23'
$ws.Range('E33').Value = '  +4.03%  '
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.90'
$ws.Range('E37').Value = '  +2.91%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '47.10'
$ws.Range('E38').Value = '  +3.68%  '
$ws.Range('E39').Value = '  +6.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '50.28'
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('E41').Value = '  +5.83%  '
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.68'
$ws.Range('E43').Value = '  +2.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('E45').Value = '  +2.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '382.67'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '2.780.37'
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '135.22'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.95'
$ws.Range('E50').Value = '  +6.41%  '
$ws.Range('E51').Value = '  +1.83%  '
